$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new sheet named "Sheet1" right after "Tabelle1"
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Sheet1"

$newSheet.Range("A1").Value = "final_cluster"
$newSheet.Range("B1").Value = "cluster_name"
$newSheet.Range("C1").Value = "Fitness"
$newSheet.Range("D1").Value = "Precision"

$data = @(
    @("01_01_01", "3-way match, invoice after GR (with SRM; Item Type: Service)", 77.3, 60.5),
    @("01_01_02", "3-way match, invoice after GR (without SRM, Item Type: Standard)", 95.8, 100),
    @("01_02_01", "3-way match, invoice after GR (without SRM; Item Type: Service)", 88.6, 80.6),
    @("01_02_02", "3-way match, invoice after GR (without SRM, Item Type: Standard)", 95, 72.7),
    @("01_02_03", "3-way match, invoice after GR (without SRM, Item Type: Subcontracting and Third-Party)", 95.4, 90.5),
    @("02_01", "3-way match, invoice before GR (with SRM)", 91.3, 99.8),
    @("02_02_01", "3-way match, invoice before GR (without SRM, Item Type: Standard)", 97.9, 100),
    @("02_02_02", "3-way match, invoice before GR (without SRM, Item Type: Subcontracting)", 90.3, 87.8),
    @("02_02_03", "3-way match, invoice before GR (without SRM, Item Type: Third-Party)", 91.6, 90.9),
    @("03", "2-way match", 97, 70.5),
    @("04", "Consignment", 100, 100)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$newSheet.Columns.Item(2).ColumnWidth = 70.89453125

# Bold header with border
$headerRange = $newSheet.Range("C1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.Item(1).LineStyle = 1
$headerRange.Borders.Item(1).Weight = -4138

$newSheet.Range("A1:D1").Borders.Item(9).LineStyle = 1
$newSheet.Range("A1:D1").Borders.Item(10).LineStyle = 1

$newSheet.Select()
